# Commit message: "corrected the spelling for non functional heading"
#
# On the Agenda slide (slide 3), the bullet that should read
# "Non Functional Requirements" was mis-typed as "Functional Requirements".
# Fix the spelling by splitting the run so the leading "Non Functional "
# is prefixed onto the existing "Requirements" text, matching the
# corrected heading while preserving the rest of the bullet formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# "Google Shape;292;p15" is the Agenda bullet-list placeholder; the 4th
# paragraph (after "Overall System Description", "Functional Requirements",
# "Use Cases") is the mis-spelled "Functional Requirements" bullet.
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(4)

# Replace the leading "Functional " (11 characters) with "Non Functional "
# so the bullet text becomes "Non Functional Requirements".
$chars = $para.Characters(1, 11)
$chars.Text = "Non Functional "
